$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sec invoice Master")

# The workbook is used as a running log: each new overage line item is
# pasted into row 2 (FC Order ID / Tracking# / New Invoice Amount /
# SECONDARY INV#), one batch of data at a time. This commit adds the
# batch of new secondary-invoice lines below, finishing with FC Order ID
# 59090385, which is what is left sitting in row 2 of "Sec invoice
# Master" afterwards.

$newRows = @(
    @("59088085", "999U967432", "59088085+1", "235.52"),
    @("59090304", "CEV1073207", "59090304+1", "274.21"),
    @("59090305", "CEV1073208", "59090305+1", "715.8"),
    @("59090308", "CEV1073209", "59090308+1", "45.54"),
    @("59090309", "CEV1073210", "59090309+1", "491.7"),
    @("59090310", "CEV1073211", "59090310+1", "713.46"),
    @("59090311", "CEV1073212", "59090311+1", "179.16"),
    @("59090312", "CEV1073213", "59090312+1", "167.98"),
    @("59090384", "CEV1073216", "59090384+1", "400.48"),
    @("59090385", "CEV1073217", "59090385+1", "520.9")
)

foreach ($row in $newRows) {
    $ws.Range("A2").Value = $row[0]
    $ws.Range("C2").Value = $row[1]
    $ws.Range("I2").Value = $row[2]
    $ws.Range("F2").Value = $row[3]
}

$ws.Columns.Item(3).ColumnWidth = 11.5
